# Horarios actualizados Línea 141 - 478
# Update the "última actualización" timestamp, row counts, and schedule
# rows on the three sheets (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

$newTime = "03:18:26"

# ---------------------------------------------------------------------
# Sheet "LP1912": schedule shifts up one stop, a new stop is appended,
# and the grid grows from A1:E9 to A1:E10.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 5"

$rows1 = @(
    @($newTime, "03:48", "14_ABASTO", 30, "LP1912"),
    @($newTime, "04:01", "81_EL PELIGRO", 43, "LP1912"),
    @($newTime, "04:46", "215A_EL PATO", 88, "LP1912"),
    @($newTime, "04:53", "11_ETCHEVERRY", 95, "LP1912"),
    @($newTime, "05:16", "17_ROMERO", 118, "LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215": only the 215A_EL PATO stop remains, updated with a
# new arrival time/minutes, the 215_ALUAR row is gone, and the grid
# shrinks from A1:E7 to A1:E6.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Cells.Item(6, 1).Value = $newTime
$ws2.Cells.Item(6, 2).Value = "04:46"
$ws2.Cells.Item(6, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(6, 4).Value = 88
$ws2.Cells.Item(6, 5).Value = "LP1912"

# Remove the now-obsolete row 7 (previously 215A_EL PATO / 113 minutes).
$ws2.Rows.Item(7).Delete()

# ---------------------------------------------------------------------
# Sheet "6203-6173": only the timestamp refreshes; still no data rows.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
